$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# 1. Delete the 12 data rows (17-28) that belong to the removed workers
#    (ORIMAR CASTILLA PAJOY x5, JAINER ARELLANES BLANQUICETT x6, MIGUEL ALBERTO ALEMAN GUILLARTE x1).
#    The remaining data row for CAMILO ANDRES MARTINEZ PEREIRA (old row 29) shifts up to row 17,
#    and the footer rows (old 34/35) shift up to rows 22/23.
$ws.Range("A17:J28").EntireRow.Delete() | Out-Null

# 2. Update the "Salario Basico" value for the remaining worker row
$ws.Range("G17").Value = 1723500

# 3. Update the summary figures
$ws.Range("E11").Value = 50823
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 2

# 4. Column D ("Nombre Trabajador") best-fit width shrinks now that the longest
#    name left in the table is shorter.
$ws.Columns.Item(4).ColumnWidth = 34
